$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3227.8
$ws.Range("I62").Value = 3512.8125
$ws.Range("J62").Value = 2087.75
$ws.Range("K62").Value = 3512.8125
$ws.Range("L62").Value = 2087.75
$ws.Range("M62").Value = -2888.8125
$ws.Range("N62").Value = -3335.75
$ws.Range("H65").Value = 3227.8
$ws.Range("I65").Value = 3512.8125
$ws.Range("J65").Value = 2087.75
$ws.Range("K65").Value = 17564.0625
$ws.Range("L65").Value = 10438.75
$ws.Range("M65").Value = -14444.0625
$ws.Range("N65").Value = -16678.75
$ws.Range("H98").Value = 41631.89
$ws.Range("J98").Value = 334402
$ws.Range("L98").Value = 334402
$ws.Range("N98").Value = -337398
$ws.Range("H122").Value = 41631.89
$ws.Range("J122").Value = 334402
$ws.Range("L122").Value = 1003206
$ws.Range("N122").Value = -1008106
$ws.Range("H128").Value = 34333.332
$ws.Range("J128").Value = 34333.332
$ws.Range("L128").Value = 34333.332
$ws.Range("N128").Value = -44293.332
$ws.Range("H132").Value = 3041639
$ws.Range("I132").Value = 3665154
$ws.Range("J132").Value = 2003.625
$ws.Range("K132").Value = 10995462
$ws.Range("L132").Value = 6010.875
$ws.Range("M132").Value = -10992932
$ws.Range("N132").Value = -11070.875
$ws.Range("H135").Value = 556.7917
$ws.Range("I135").Value = 504.8085
$ws.Range("K135").Value = 4543.2765
$ws.Range("M135").Value = -2008.2765
$ws.Range("H137").Value = 1106.4678
$ws.Range("I137").Value = 868.60785
$ws.Range("K137").Value = 2605.82355
$ws.Range("M137").Value = -55.82355000000007
$ws.Range("H138").Value = 1726.065
$ws.Range("I138").Value = 1318.1892
$ws.Range("J138").Value = 2103.35
$ws.Range("K138").Value = 3954.5676
$ws.Range("L138").Value = 6310.049999999999
$ws.Range("M138").Value = 1185.4324
$ws.Range("N138").Value = -16590.05
$ws.Range("H141").Value = 1829.7963
$ws.Range("I141").Value = 1143.1515
$ws.Range("J141").Value = 2908.8096
$ws.Range("K141").Value = 3429.4545
$ws.Range("L141").Value = 8726.4288
$ws.Range("M141").Value = 1750.5455
$ws.Range("N141").Value = -19086.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1141.6364
$ws.Range("I61").Value = 694.2381
$ws.Range("J61").Value = 2587.077
$ws.Range("K61").Value = 694.2381
$ws.Range("L61").Value = 2587.077
$ws.Range("M61").Value = -482.2381
$ws.Range("N61").Value = -3011.077
$ws.Range("H74").Value = 881.6981
$ws.Range("I74").Value = 839.375
$ws.Range("J74").Value = 1011.9231
$ws.Range("K74").Value = 839.375
$ws.Range("L74").Value = 1011.9231
$ws.Range("M74").Value = 34.625
$ws.Range("N74").Value = -2759.9231
$ws.Range("H77").Value = 881.6981
$ws.Range("I77").Value = 839.375
$ws.Range("J77").Value = 1011.9231
$ws.Range("K77").Value = 4196.875
$ws.Range("L77").Value = 5059.6155
$ws.Range("M77").Value = 171.125
$ws.Range("N77").Value = -13795.6155
$ws.Range("H94").Value = 34915
$ws.Range("J94").Value = 34915
$ws.Range("L94").Value = 34915
$ws.Range("N94").Value = -36717
$ws.Range("H132").Value = 10529.308
$ws.Range("I132").Value = 12228.1
$ws.Range("J132").Value = 4866.6665
$ws.Range("K132").Value = 36684.3
$ws.Range("L132").Value = 14599.9995
$ws.Range("M132").Value = -34154.3
$ws.Range("N132").Value = -19659.9995
$ws.Range("H136").Value = 1141.6364
$ws.Range("I136").Value = 694.2381
$ws.Range("J136").Value = 2587.077
$ws.Range("K136").Value = 2082.7143
$ws.Range("L136").Value = 7761.231000000001
$ws.Range("M136").Value = 467.2856999999999
$ws.Range("N136").Value = -12861.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 2020
$ws.Range("I128").Value = 2020
$ws.Range("K128").Value = 6060
$ws.Range("M128").Value = -3570
$ws.Range("H134").Value = 25266.117
$ws.Range("I134").Value = 35158.3
$ws.Range("J134").Value = 2438
$ws.Range("K134").Value = 105474.9
$ws.Range("L134").Value = 7314
$ws.Range("M134").Value = -102939.9
$ws.Range("N134").Value = -12384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H23").Value = 4000
$ws.Range("J23").Value = 4000
$ws.Range("L23").Value = 4000
$ws.Range("N23").Value = -4480
$ws.Range("H27").Value = 4000
$ws.Range("J27").Value = 4000
$ws.Range("L27").Value = 4000
$ws.Range("N27").Value = -4384
$ws.Range("H31").Value = 5558200.5
$ws.Range("I31").Value = 1813.909
$ws.Range("J31").Value = 20838264
$ws.Range("K31").Value = 1813.909
$ws.Range("L31").Value = 20838264
$ws.Range("M31").Value = -1518.909
$ws.Range("N31").Value = -20838854
$ws.Range("H34").Value = 5558200.5
$ws.Range("I34").Value = 1813.909
$ws.Range("J34").Value = 20838264
$ws.Range("K34").Value = 1813.909
$ws.Range("L34").Value = 20838264
$ws.Range("M34").Value = -1611.909
$ws.Range("N34").Value = -20838668
$ws.Range("H58").Value = 708.1389
$ws.Range("I58").Value = 749.63635
$ws.Range("J58").Value = 642.9286
$ws.Range("K58").Value = 749.63635
$ws.Range("L58").Value = 642.9286
$ws.Range("M58").Value = -546.63635
$ws.Range("N58").Value = -1048.9286
$ws.Range("H102").Value = 38300
$ws.Range("J102").Value = 38300
$ws.Range("L102").Value = 38300
$ws.Range("N102").Value = -43168
$ws.Range("H132").Value = 1458.3043
$ws.Range("I132").Value = 1410.898
$ws.Range("J132").Value = 1574.45
$ws.Range("K132").Value = 4232.694
$ws.Range("L132").Value = 4723.35
$ws.Range("M132").Value = -1702.694
$ws.Range("N132").Value = -9783.35
$ws.Range("H134").Value = 931.16
$ws.Range("I134").Value = 885.3182
$ws.Range("J134").Value = 1267.3334
$ws.Range("K134").Value = 2655.9546
$ws.Range("L134").Value = 3802.0002
$ws.Range("M134").Value = -120.9546
$ws.Range("N134").Value = -8872.0002
$ws.Range("H136").Value = 708.1389
$ws.Range("I136").Value = 749.63635
$ws.Range("J136").Value = 642.9286
$ws.Range("K136").Value = 2248.90905
$ws.Range("L136").Value = 1928.7858
$ws.Range("M136").Value = 301.0909499999998
$ws.Range("N136").Value = -7028.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 674.75
$ws.Range("I15").Value = 200
$ws.Range("J15").Value = 833
$ws.Range("K15").Value = 600
$ws.Range("L15").Value = 2499
$ws.Range("M15").Value = -460
$ws.Range("N15").Value = -2779
$ws.Range("H20").Value = 1000
$ws.Range("J20").Value = 1000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3454
$ws.Range("H122").Value = 1477.625
$ws.Range("I122").Value = 1729.1428
$ws.Range("J122").Value = 1282
$ws.Range("K122").Value = 15562.2852
$ws.Range("L122").Value = 11538
$ws.Range("M122").Value = -13112.2852
$ws.Range("N122").Value = -16438

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 8660
$ws.Range("J53").Value = 9575
$ws.Range("L53").Value = 9575
$ws.Range("N53").Value = -10837
$ws.Range("H58").Value = 20000000
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H132").Value = 24184.021
$ws.Range("I132").Value = 32117.787
$ws.Range("J132").Value = 2366.1667
$ws.Range("K132").Value = 96353.361
$ws.Range("L132").Value = 7098.500100000001
$ws.Range("M132").Value = -93823.361
$ws.Range("N132").Value = -12158.5001
$ws.Range("H134").Value = 23950
$ws.Range("J134").Value = 23950
$ws.Range("L134").Value = 71850
$ws.Range("N134").Value = -76920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3400.5715
$ws.Range("I40").Value = 2476
$ws.Range("J40").Value = 4633.3335
$ws.Range("K40").Value = 2476
$ws.Range("L40").Value = 4633.3335
$ws.Range("M40").Value = -2340
$ws.Range("N40").Value = -4905.3335
$ws.Range("H55").Value = 277.95456
$ws.Range("I55").Value = 287.27274
$ws.Range("J55").Value = 268.63635
$ws.Range("K55").Value = 287.27274
$ws.Range("L55").Value = 268.63635
$ws.Range("M55").Value = -114.27274
$ws.Range("N55").Value = -614.63635
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H104").Value = 29929.4
$ws.Range("J104").Value = 29929.4
$ws.Range("L104").Value = 29929.4
$ws.Range("N104").Value = -36917.4
$ws.Range("H122").Value = 2970.1738
$ws.Range("I122").Value = 3062.2307
$ws.Range("J122").Value = 2850.5
$ws.Range("K122").Value = 9186.6921
$ws.Range("L122").Value = 8551.5
$ws.Range("M122").Value = -6736.6921
$ws.Range("N122").Value = -13451.5
$ws.Range("H132").Value = 4957.418
$ws.Range("I132").Value = 6707.0835
$ws.Range("J132").Value = 1642.2632
$ws.Range("K132").Value = 20121.2505
$ws.Range("L132").Value = 4926.7896
$ws.Range("M132").Value = -17591.2505
$ws.Range("N132").Value = -9986.7896
$ws.Range("H135").Value = 45420
$ws.Range("J135").Value = 45420
$ws.Range("L135").Value = 45420
$ws.Range("N135").Value = -55560
$ws.Range("H136").Value = 3058.8164
$ws.Range("I136").Value = 3294.3784
$ws.Range("J136").Value = 2332.5
$ws.Range("K136").Value = 9883.135200000001
$ws.Range("L136").Value = 6997.5
$ws.Range("M136").Value = -7333.135200000001
$ws.Range("N136").Value = -12097.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4074.75
$ws.Range("I81").Value = 5033
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 10066
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -9005
$ws.Range("N81").Value = -4522
$ws.Range("H84").Value = 4074.75
$ws.Range("I84").Value = 5033
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 50330
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -45026
$ws.Range("N84").Value = -22608
$ws.Range("H136").Value = 3177.132
$ws.Range("I136").Value = 3473.818
$ws.Range("J136").Value = 1726.6666
$ws.Range("K136").Value = 10421.454
$ws.Range("L136").Value = 5179.9998
$ws.Range("M136").Value = -7871.454000000002
$ws.Range("N136").Value = -10279.9998
